$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 137, shifting existing rows 137:142 down to 138:143.
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row 137 with the new weekly record.
$ws.Range("A137").Value = 11
$ws.Range("B137").Value = "Vega Monumental Concepción"
$ws.Range("C137").Value = "Bíobío"
$ws.Range("D137").Value = 44610
$ws.Range("E137").Value = 8
$ws.Range("F137").Value = 100112003
$ws.Range("G137").Value = "Ajo"
$ws.Range("H137").Value = "Chino"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 400
$ws.Range("K137").Value = 17000
$ws.Range("L137").Value = 18000
$ws.Range("M137").Value = 17500
$ws.Range("N137").Value = "`$/caja 10 kilos"
$ws.Range("O137").Value = "China"
$ws.Range("P137").Value = 1750
$ws.Range("Q137").Value = 10
$ws.Range("R137").Value = "Hortaliza"
